$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1107.6923
$ws.Range("I18").Value = 1036.3636
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 1036.3636
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = -752.3635999999999
$ws.Range("N18").Value = -2068
$ws.Range("H100").Value = 26316688
$ws.Range("I100").Value = 960.38464
$ws.Range("J100").Value = 83334100
$ws.Range("K100").Value = 960.38464
$ws.Range("L100").Value = 83334100
$ws.Range("M100").Value = -419.38464
$ws.Range("N100").Value = -83335182
$ws.Range("H107").Value = 919.25
$ws.Range("I107").Value = 953.5
$ws.Range("J107").Value = 850.75
$ws.Range("K107").Value = 953.5
$ws.Range("L107").Value = 850.75
$ws.Range("M107").Value = 966.5
$ws.Range("N107").Value = -4690.75
$ws.Range("H111").Value = 2857.2104
$ws.Range("I111").Value = 2299
$ws.Range("J111").Value = 4066.6667
$ws.Range("K111").Value = 6897
$ws.Range("L111").Value = 12200.0001
$ws.Range("M111").Value = -3830
$ws.Range("N111").Value = -18334.0001
$ws.Range("H138").Value = 2664.59
$ws.Range("I138").Value = 1146.0526
$ws.Range("J138").Value = 3595.3064
$ws.Range("K138").Value = 3438.1578
$ws.Range("L138").Value = 10785.9192
$ws.Range("M138").Value = 1701.8422
$ws.Range("N138").Value = -21065.9192
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22637.86
$ws.Range("I32").Value = 8791.279
$ws.Range("J32").Value = 114238.305
$ws.Range("K32").Value = 8791.279
$ws.Range("L32").Value = 114238.305
$ws.Range("M32").Value = -8504.279
$ws.Range("N32").Value = -114812.305
$ws.Range("H44").Value = 27899.5
$ws.Range("H97").Value = 738.25714
$ws.Range("I97").Value = 454.08334
$ws.Range("K97").Value = 454.08334
$ws.Range("M97").Value = 41.91665999999998
$ws.Range("H107").Value = 23742.666
$ws.Range("J107").Value = 23742.666
$ws.Range("L107").Value = 23742.666
$ws.Range("N107").Value = -31422.666
$ws.Range("H123").Value = 656875
$ws.Range("J123").Value = 656875
$ws.Range("L123").Value = 656875
$ws.Range("N123").Value = -666675
$ws.Range("H133").Value = 44347.75
$ws.Range("J133").Value = 44347.75
$ws.Range("L133").Value = 44347.75
$ws.Range("N133").Value = -49407.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4781.8335
$ws.Range("I20").Value = 6691.857
$ws.Range("J20").Value = 2107.8
$ws.Range("K20").Value = 6691.857
$ws.Range("L20").Value = 2107.8
$ws.Range("M20").Value = -6444.857
$ws.Range("N20").Value = -2601.8
$ws.Range("H35").Value = 35037
$ws.Range("J35").Value = 35037
$ws.Range("L35").Value = 35037
$ws.Range("N35").Value = -35657
$ws.Range("H82").Value = 17182.445
$ws.Range("I82").Value = 4038.2727
$ws.Range("J82").Value = 37837.57
$ws.Range("K82").Value = 4038.2727
$ws.Range("L82").Value = 37837.57
$ws.Range("M82").Value = -3655.2727
$ws.Range("N82").Value = -38603.57
$ws.Range("H85").Value = 17182.445
$ws.Range("I85").Value = 4038.2727
$ws.Range("J85").Value = 37837.57
$ws.Range("K85").Value = 4038.2727
$ws.Range("L85").Value = 37837.57
$ws.Range("M85").Value = -2712.2727
$ws.Range("N85").Value = -40489.57
$ws.Range("H122").Value = 42500
$ws.Range("J122").Value = 42500
$ws.Range("L122").Value = 42500
$ws.Range("N122").Value = -52300
$ws.Range("H125").Value = 43886.668
$ws.Range("J125").Value = 43886.668
$ws.Range("L125").Value = 43886.668
$ws.Range("N125").Value = -53726.668
$ws.Range("H126").Value = 27998
$ws.Range("J126").Value = 27998
$ws.Range("L126").Value = 27998
$ws.Range("N126").Value = -37878
$ws.Range("H132").Value = 39939.9
$ws.Range("J132").Value = 39939.9
$ws.Range("L132").Value = 39939.9
$ws.Range("N132").Value = -50059.9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3332.0598
$ws.Range("I31").Value = 1612.6923
$ws.Range("J31").Value = 5726.893
$ws.Range("K31").Value = 1612.6923
$ws.Range("L31").Value = 5726.893
$ws.Range("M31").Value = -1317.6923
$ws.Range("N31").Value = -6316.893
$ws.Range("H34").Value = 3332.0598
$ws.Range("I34").Value = 1612.6923
$ws.Range("J34").Value = 5726.893
$ws.Range("K34").Value = 1612.6923
$ws.Range("L34").Value = 5726.893
$ws.Range("M34").Value = -1410.6923
$ws.Range("N34").Value = -6130.893
$ws.Range("H41").Value = 16280.75
$ws.Range("J41").Value = 25032.5
$ws.Range("L41").Value = 25032.5
$ws.Range("N41").Value = -25888.5
$ws.Range("H50").Value = 9213.6
$ws.Range("J50").Value = 9213.6
$ws.Range("L50").Value = 9213.6
$ws.Range("N50").Value = -10463.6
$ws.Range("H51").Value = 9306.5
$ws.Range("J51").Value = 9306.5
$ws.Range("L51").Value = 9306.5
$ws.Range("N51").Value = -10778.5
$ws.Range("H59").Value = 15097.3
$ws.Range("I59").Value = 12000
$ws.Range("J59").Value = 15871.625
$ws.Range("K59").Value = 12000
$ws.Range("L59").Value = 15871.625
$ws.Range("M59").Value = -10855
$ws.Range("N59").Value = -18161.625
$ws.Range("H60").Value = 24931.23
$ws.Range("J60").Value = 24931.23
$ws.Range("L60").Value = 24931.23
$ws.Range("N60").Value = -25953.23
$ws.Range("H61").Value = 9306.5
$ws.Range("J61").Value = 9306.5
$ws.Range("L61").Value = 9306.5
$ws.Range("N61").Value = -10002.5
$ws.Range("H63").Value = 30750
$ws.Range("J63").Value = 30750
$ws.Range("L63").Value = 30750
$ws.Range("N63").Value = -32122
$ws.Range("H66").Value = 30750
$ws.Range("J66").Value = 30750
$ws.Range("L66").Value = 92250
$ws.Range("N66").Value = -99114
$ws.Range("H68").Value = 17999.666
$ws.Range("J68").Value = 17999.666
$ws.Range("L68").Value = 17999.666
$ws.Range("N68").Value = -19497.666
$ws.Range("H71").Value = 17999.666
$ws.Range("J71").Value = 17999.666
$ws.Range("L71").Value = 53998.99800000001
$ws.Range("N71").Value = -61486.99800000001
$ws.Range("H74").Value = 1213383.8
$ws.Range("J74").Value = 1213383.8
$ws.Range("L74").Value = 1213383.8
$ws.Range("N74").Value = -1215131.8
$ws.Range("H77").Value = 1213383.8
$ws.Range("J77").Value = 1213383.8
$ws.Range("L77").Value = 3640151.4
$ws.Range("N77").Value = -3648887.4
$ws.Range("H97").Value = 17456
$ws.Range("J97").Value = 17456
$ws.Range("L97").Value = 17456
$ws.Range("N97").Value = -19438
$ws.Range("H130").Value = 56385
$ws.Range("J130").Value = 56385
$ws.Range("L130").Value = 56385
$ws.Range("N130").Value = -66425
$ws.Range("H131").Value = 45323
$ws.Range("J131").Value = 45323
$ws.Range("L131").Value = 45323
$ws.Range("N131").Value = -55403
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2219.2727
$ws.Range("I137").Value = 1487.8948
$ws.Range("J137").Value = 3211.8572
$ws.Range("K137").Value = 4463.6844
$ws.Range("L137").Value = 9635.571599999999
$ws.Range("M137").Value = 636.3155999999999
$ws.Range("N137").Value = -19835.5716
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14798.5
$ws.Range("J57").Value = 17164.666
$ws.Range("L57").Value = 17164.666
$ws.Range("N57").Value = -18804.666
$ws.Range("H70").Value = 34755.883
$ws.Range("I70").Value = 45752.16
$ws.Range("J70").Value = 4210.6665
$ws.Range("K70").Value = 45752.16
$ws.Range("L70").Value = 4210.6665
$ws.Range("M70").Value = -45482.16
$ws.Range("N70").Value = -4750.6665
$ws.Range("H73").Value = 34755.883
$ws.Range("I73").Value = 45752.16
$ws.Range("J73").Value = 4210.6665
$ws.Range("K73").Value = 45752.16
$ws.Range("L73").Value = 4210.6665
$ws.Range("M73").Value = -44816.16
$ws.Range("N73").Value = -6082.6665
$ws.Range("H123").Value = 15133.917
$ws.Range("J123").Value = 15133.917
$ws.Range("L123").Value = 15133.917
$ws.Range("N123").Value = -20033.917
$ws.Range("H124").Value = 42250
$ws.Range("J124").Value = 42250
$ws.Range("L124").Value = 42250
$ws.Range("N124").Value = -52070
$ws.Range("H128").Value = 46710.25
$ws.Range("J128").Value = 46710.25
$ws.Range("L128").Value = 46710.25
$ws.Range("N128").Value = -56670.25
$ws.Range("H130").Value = 49056
$ws.Range("J130").Value = 49056
$ws.Range("L130").Value = 49056
$ws.Range("N130").Value = -59096
$ws.Range("H133").Value = 44624.445
$ws.Range("J133").Value = 44624.445
$ws.Range("L133").Value = 44624.445
$ws.Range("N133").Value = -54744.445
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 20477.6
$ws.Range("J96").Value = 20477.6
$ws.Range("L96").Value = 20477.6
$ws.Range("N96").Value = -25969.6
$ws.Range("H122").Value = 4205.1665
$ws.Range("I122").Value = 3644
$ws.Range("J122").Value = 4562.273
$ws.Range("K122").Value = 10932
$ws.Range("L122").Value = 13686.819
$ws.Range("M122").Value = -8482
$ws.Range("N122").Value = -18586.819
$ws.Range("H132").Value = 3705.2046
$ws.Range("I132").Value = 4039.4075
$ws.Range("J132").Value = 3174.4119
$ws.Range("K132").Value = 12118.2225
$ws.Range("L132").Value = 9523.235700000001
$ws.Range("M132").Value = -9588.2225
$ws.Range("N132").Value = -14583.2357
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 19333.334
$ws.Range("J64").Value = 19333.334
$ws.Range("L64").Value = 19333.334
$ws.Range("N64").Value = -19829.334
$ws.Range("H67").Value = 19333.334
$ws.Range("J67").Value = 19333.334
$ws.Range("L67").Value = 19333.334
$ws.Range("N67").Value = -21049.334
$ws.Range("H125").Value = 27997.273
$ws.Range("J125").Value = 27997.273
$ws.Range("L125").Value = 27997.273
$ws.Range("N125").Value = -37837.273
$ws.Range("H127").Value = 20198.428
$ws.Range("J127").Value = 20198.428
$ws.Range("L127").Value = 20198.428
$ws.Range("N127").Value = -30118.428

Write-Output "Applied 256 cell updates across 8 sheets"